$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$ws.Range("B21").Value = "5-75 Manufacturing " + [char]10 + "5-30 Services and others"
$ws.Range("D21").Value = "300,000-15Million RM Manufacturing " + [char]10 + "300,000-3Million RM Services & others"

$ws.Range("B22").Value = "75-200 Manufacturing " + [char]10 + "30-75 Services  and others"
$ws.Range("D22").Value = "15-50Million RM Manufacturing " + [char]10 + "3-20Million RM Services & others"

$ws.Range("B23").Value = ">200 Manufacturing " + [char]10 + ">75 Services  and others"
$ws.Range("D23").Value = ">50Million RM Manufacturing " + [char]10 + ">20Million RM Services & others"
